$d = $word.ActiveDocument

# 1. Remove the " see appendix 1" run entirely from the "OLS Regression:" paragraph.
$rng1 = $d.Content
$rng1.Find.Execute(" see appendix 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng1.Delete()

# 2. Replace " see appendix 2" with a single space, preserving the underline
#    formatting, while keeping it as its own run (not merged with the
#    preceding ":" run).
$rng2 = $d.Content
$rng2.Find.Execute(" see appendix 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng2.Delete()
$newRng = $rng2.Duplicate
$newRng.InsertAfter(" ")
$newRng.Font.Underline = 1
